$d = $word.ActiveDocument

# Locate the last bullet ("Can we make a game out of Bayesian
# inference?") -- that's where the new list item goes, right after it.
$target = $null
$targetIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*Can we make a game out of Bayesian inference?*") {
        $target = $p
        $targetIndex = $i
    }
}

# Paragraph.Range.Text includes the trailing paragraph-mark (chr 13);
# strip it so Find/Replace matches only the visible sentence.
$bulletText = $target.Range.Text.Substring(0, $target.Range.Text.Length - 1)

# Split the paragraph right after its text (before the trailing pilcrow)
# so a brand-new, empty "ListParagraph" bullet paragraph follows it --
# same pPr (style/numbering) gets carried over automatically, just like
# pressing Enter at the end of the bullet in the UI.
$found = $d.Content.Find.Execute($bulletText, $true, $false, $false, `
    $false, $false, $true, 1, $false, ($bulletText + "^p"), 2)

# Fill in the new (now-empty) paragraph that was just created.
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = "Conjugate prior prior, likelihood, posterior, posterior predictive tool."
